# Added Jumping on Clouds and *linked to educative project
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Path drawing for the "Jumping on Clouds" walk (rows 11-14) ---
# Leading apostrophes force Excel's quote-prefix text entry for the
# "/" characters, matching how the existing Counting-Valleys drawing
# (rows 4-6) stores its slashes.
$ws.Range("C11").Value = "'/"
$ws.Range("D11").Value = "\"

$ws.Range("A12").Value = "_"
$ws.Range("B12").Value = "'/"
$ws.Range("E12").Value = "\"

$ws.Range("F13").Value = "\"
$ws.Range("J13").Value = "'/"

$ws.Range("G14").Value = "'/"
$ws.Range("H14").Value = "\"
$ws.Range("I14").Value = "'/"

# --- Jumping on Clouds input data (row 16: U/D directions, row 17: index) ---
$ws.Range("B16").Value = "U"
$ws.Range("C16").Value = "U"
$ws.Range("D16").Value = "D"
$ws.Range("E16").Value = "D"
$ws.Range("F16").Value = "D"
$ws.Range("G16").Value = "U"
$ws.Range("H16").Value = "D"
$ws.Range("I16").Value = "U"
$ws.Range("J16").Value = "U"

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 7
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 9

# Leave the selection the same way the author left it after entering the
# new block (B16:J16, active cell B16).
[void]$ws.Range("B16:J16").Select()
